$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '27.601.57'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.09'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '315.36'
$ws.Range("E5").Value = '  +0.75%  '

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4312'
$ws.Range("E7").Value = '  +1.01%  '

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3693'
$ws.Range("E8").Value = '  +1.90%  '

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.07320'
$ws.Range("E9").Value = '  +0.28%  '

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '0.8764'
$ws.Range("E10").Value = '  +0.21%  '

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '21.01'
$ws.Range("E11").Value = '  +2.05%  '

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.917.71'
$ws.Range("E12").Value = '  +1.52%  '

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.474'
$ws.Range("E13").Value = '  +2.91%  '

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '6.603'
$ws.Range("E14").Value = '  +1.52%  '

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = '0.06948'
$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '81.33'
$ws.Range("E17").Value = '  +1.67%  '

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.000009058'
$ws.Range("E18").Value = '  +0.61%  '

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '15.60'
$ws.Range("E20").Value = '  +1.63%  '

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '28.165.08'
$ws.Range("E21").Value = '  +1.76%  '

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.090'
$ws.Range("E22").Value = '  +2.75%  '

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '10.95'
$ws.Range("E23").Value = '  +5.53%  '

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.170.53'
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '1.989'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '154.10'
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.99'
$ws.Range("E27").Value = '  +1.00%  '

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '5.312'
$ws.Range("E28").Value = '  +1.17%  '

$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '115.92'
$ws.Range("E29").Value = '  -4.68%  '

$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '1.878'
$ws.Range("E30").Value = '  +1.35%  '

$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.08905'
$ws.Range("E31").Value = '  -0.19%  '

$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.7876'
$ws.Range("E32").Value = '  +3.73%  '

$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.608'
$ws.Range("E33").Value = '  +1.94%  '

$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.984'
$ws.Range("E34").Value = '  +0.48%  '

$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.169'
$ws.Range("E35").Value = '  +6.53%  '

$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.109'
$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05439'
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01966'
$ws.Range("E39").Value = '  +1.77%  '

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.845'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.5175'
$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1694'
$ws.Range("E42").Value = '  +2.54%  '

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.789'
$ws.Range("E43").Value = '  +0.51%  '

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '8.665'
$ws.Range("E44").Value = '  +4.08%  '

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '10.72'
$ws.Range("E45").Value = '  +4.17%  '

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.4798'
$ws.Range("E46").Value = '  +2.89%  '

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '106.62'
$ws.Range("E47").Value = '  +1.61%  '

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06538'
$ws.Range("E48").Value = '  -0.27%  '

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").Value = '1.002'
$ws.Range("E49").Value = '  -0.10%  '

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.664'
$ws.Range("E50").Value = '  +3.03%  '

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.847'
$ws.Range("E51").Value = '  +6.31%  '
